$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Before A Proper Officer at <<courtName>>", $true, $false, $false, $false, $false,
    $true, 1, $false, "Before A Proper Officer at <<siteName>> - <<address>> - <<postcode>>.", 2
)
